# Scrum Sprint Planning - update the "Documentatie" sprint row:
#   - Fix/rename the task name from "Documentatie" to "Documentation"
#   - Mark it as done ("n" -> "y")
#   - Leave the final selection on B3 (matches the saved cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Done" flag for the Documentation row first so the new "y" string
# is registered in the shared-strings table before the renamed task text,
# matching the order the strings were appended to the table.
$ws.Range("C2").Value = "y"
$ws.Range("B2").Value = "Documentation"

# Leave the selection where the edit finished, on the renamed cell's row.
$ws.Range("B3").Select() | Out-Null
